# Update row 6 (first data row) of the "Ds cấp phát TSCĐ" sheet with a new
# allocation record: date, asset code, warehouse code, and receiving
# department code.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ds cấp phát TSCĐ")

$ws.Range("B6").Value = "20/09/2024"
$ws.Range("C6").Value = "TS-007855"
$ws.Range("D6").Value = "K.CCDC"
$ws.Range("E6").Value = "A1"
